$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation arrived for this market/category. It belongs
# chronologically at row 45 (28-Mar-2023), so insert a fresh row there and
# push the existing row 45 (18-Feb-2022) down to row 46, preserving its data.
$ws.Rows.Item(45).Insert()

# Fill the newly inserted row 45 with the new weekly observation.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 45013
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100114007
$ws.Range("G45").Value = "Jengibre"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 220
$ws.Range("K45").Value = 15000
$ws.Range("L45").Value = 16000
$ws.Range("M45").Value = 15455
$ws.Range("N45").Value = "$/caja 13 kilos"
$ws.Range("O45").Value = "Perú"
$ws.Range("P45").Value = 1189
$ws.Range("Q45").Value = 13
$ws.Range("R45").Value = "Hortaliza"
